# Updated cryptos list on Wed Aug  7 15:49:54 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-style pattern used below for cells whose new text would otherwise be
# auto-converted to a number by Excel: force Text format, assign the value,
# then clear the (temporary) formatting so the cell keeps looking like the
# untouched ones (no explicit style index).

# Rows 42 and 43 swapped coins (Hedera and Filecoin traded places in the
# ranking) with updated price/volume figures.
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.38"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.18%  "

$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0542"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.71%  "

# Updated price and volume figures for the remaining rows.
$ws.Range("D2").Value = "55.863.35"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "2.392.02"
$ws.Range("E3").Value = "  -4.10%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "479.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.42"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.81%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E8").Value = "  -1.73%  "
$ws.Range("D9").Value = "2.390.74"
$ws.Range("E9").Value = "  -4.86%  "
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("E11").Value = "  -3.74%  "
$ws.Range("E12").Value = "  -1.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.125"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "2.807.02"
$ws.Range("E14").Value = "  -4.16%  "
$ws.Range("D15").Value = "56.284.18"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.32"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.99%  "
$ws.Range("E17").Value = "  -1.84%  "
$ws.Range("D18").Value = "2.394.28"
$ws.Range("E18").Value = "  -4.42%  "
$ws.Range("E19").Value = "  +0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "315.03"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.76"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  -1.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "56.81"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.46%  "
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.395"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.38%  "
$ws.Range("E27").Value = "  -4.14%  "
$ws.Range("D28").Value = "2.497.12"
$ws.Range("E28").Value = "  -4.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("D30").Value = "0.0₃0772"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.27"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.99"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.03%  "
$ws.Range("E36").Value = "  -2.95%  "
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.42"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.28%  "
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0947"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.584"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.22"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.62"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "253.78"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.05"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.21%  "
$ws.Range("D51").Value = "1.771.14"
$ws.Range("E51").Value = "  -7.66%  "
